$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 2-6 with revised financial figures ---
# Row 2
$ws.Range("D2").Value = 823
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 16
$ws.Range("G2").Value = 11
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 8
$ws.Range("K2").Value = 979
$ws.Range("L2").Value = 353
$ws.Range("M2").Value = 627
$ws.Range("N2").Value = 627
$ws.Range("P2").Value = 85
$ws.Range("Q2").Value = 39
$ws.Range("R2").Value = -24
$ws.Range("S2").Value = -33
$ws.Range("T2").Value = 22
$ws.Range("U2").Value = 17
$ws.Range("V2").Value = 226
$ws.Range("W2").Value = 1.98
$ws.Range("X2").Value = 0.92
$ws.Range("Y2").Value = 1.21
$ws.Range("Z2").Value = 0.76
$ws.Range("AA2").Value = 56.28
$ws.Range("AB2").Value = 686.63
$ws.Range("AC2").Value = 45
$ws.Range("AD2").Value = 54.6
$ws.Range("AE2").Value = 4048
$ws.Range("AF2").Value = 0.6
$ws.Range("AG2").Value = 45
$ws.Range("AH2").Value = 1.85
$ws.Range("AI2").Value = 92.6
$ws.Range("AJ2").Value = 15680000

# Row 3
$ws.Range("D3").Value = 769
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 7
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("K3").Value = 973
$ws.Range("L3").Value = 338
$ws.Range("M3").Value = 635
$ws.Range("N3").Value = 635
$ws.Range("P3").Value = 85
$ws.Range("Q3").Value = 57
$ws.Range("R3").Value = -30
$ws.Range("S3").Value = -34
$ws.Range("T3").Value = 16
$ws.Range("U3").Value = 41
$ws.Range("V3").Value = 198
$ws.Range("W3").Value = 0.71
$ws.Range("X3").Value = 0.68
$ws.Range("Y3").Value = 0.83
$ws.Range("Z3").Value = 0.53
$ws.Range("AA3").Value = 53.2
$ws.Range("AB3").Value = 686.77
$ws.Range("AC3").Value = 31
$ws.Range("AD3").Value = 285.94
$ws.Range("AE3").Value = 4101
$ws.Range("AF3").Value = 2.13
$ws.Range("AG3").Value = 20
$ws.Range("AH3").Value = 0.23
$ws.Range("AI3").Value = 60.6
$ws.Range("AJ3").Value = 15680000

# Row 4
$ws.Range("D4").Value = 1038
$ws.Range("E4").Value = 34
$ws.Range("F4").Value = 34
$ws.Range("G4").Value = 34
$ws.Range("H4").Value = 29
$ws.Range("I4").Value = 29
$ws.Range("K4").Value = 1082
$ws.Range("L4").Value = 423
$ws.Range("M4").Value = 659
$ws.Range("N4").Value = 659
$ws.Range("P4").Value = 85
$ws.Range("Q4").Value = 33
$ws.Range("R4").Value = -69
$ws.Range("S4").Value = 54
$ws.Range("T4").Value = 56
$ws.Range("U4").Value = -24
$ws.Range("V4").Value = 255
$ws.Range("W4").Value = 3.25
$ws.Range("X4").Value = 2.83
$ws.Range("Y4").Value = 4.55
$ws.Range("Z4").Value = 2.86
$ws.Range("AA4").Value = 64.24
$ws.Range("AB4").Value = 718.47
$ws.Range("AC4").Value = 172
$ws.Range("AD4").Value = 37.14
$ws.Range("AE4").Value = 4256
$ws.Range("AF4").Value = 1.5
$ws.Range("AG4").Value = 45
$ws.Range("AH4").Value = 0.7
$ws.Range("AI4").Value = 23.92
$ws.Range("AJ4").Value = 15680000

# Row 5
$ws.Range("D5").Value = 1007
$ws.Range("E5").Value = 18
$ws.Range("F5").Value = 18
$ws.Range("G5").Value = 15
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 5
$ws.Range("K5").Value = 1054
$ws.Range("L5").Value = 401
$ws.Range("M5").Value = 653
$ws.Range("N5").Value = 653
$ws.Range("P5").Value = 85
$ws.Range("Q5").Value = 48
$ws.Range("R5").Value = -59
$ws.Range("S5").Value = 6
$ws.Range("T5").Value = 66
$ws.Range("U5").Value = -18
$ws.Range("V5").Value = 268
$ws.Range("W5").Value = 1.76
$ws.Range("X5").Value = 0.52
$ws.Range("Y5").Value = 0.8
$ws.Range("Z5").Value = 0.49
$ws.Range("AA5").Value = 61.37
$ws.Range("AB5").Value = 718.3
$ws.Range("AC5").Value = 31
$ws.Range("AD5").Value = 120.26
$ws.Range("AE5").Value = 4221
$ws.Range("AF5").Value = 0.87
$ws.Range("AG5").Value = 20
$ws.Range("AH5").Value = 0.54
$ws.Range("AI5").Value = 60.52
$ws.Range("AJ5").Value = 15680000

# Row 6
$ws.Range("D6").Value = 837
$ws.Range("E6").Value = -23
$ws.Range("F6").Value = -23
$ws.Range("G6").Value = -24
$ws.Range("H6").Value = -20
$ws.Range("I6").Value = -20
$ws.Range("K6").Value = 1006
$ws.Range("L6").Value = 377
$ws.Range("M6").Value = 628
$ws.Range("N6").Value = 628
$ws.Range("P6").Value = 85
$ws.Range("Q6").Value = -5
$ws.Range("R6").Value = -29
$ws.Range("S6").Value = -1
$ws.Range("T6").Value = 30
$ws.Range("U6").Value = -34
$ws.Range("V6").Value = 270
$ws.Range("W6").Value = -2.71
$ws.Range("X6").Value = -2.41
$ws.Range("Y6").Value = -3.16
$ws.Range("Z6").Value = -1.96
$ws.Range("AA6").Value = 60.01
$ws.Range("AB6").Value = 693.26
$ws.Range("AC6").Value = -118
$ws.Range("AD6").Value = -29.55
$ws.Range("AE6").Value = 4059
$ws.Range("AF6").Value = 0.86
$ws.Range("AG6").Value = 20
$ws.Range("AH6").Value = 0.57
$ws.Range("AI6").Value = -15.66
$ws.Range("AJ6").Value = 15680000

# --- Drop the now-unused 당기순이익(비지배)/자본총계(비지배) figures (cols J, O) for rows 2-5 ---
$ws.Range("J2:J5").ClearContents()
$ws.Range("O2:O5").ClearContents()

# --- Remove the stale forecast rows (2019E-2021E) data, keeping only the label columns A:C ---
$ws.Range("D7:AJ9").ClearContents()